# Polish slides with step-by-step adherence to outline and precise metrics
$p = $ppt.ActivePresentation
$CR = [char]13

# ---------------------------------------------------------------------
# Slide 2: Business Context & Objectives
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(6).Runs(1).Text = "Identify precursor signals of cancellation/success early."

# ---------------------------------------------------------------------
# Slide 3: Data Overview
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Paragraphs(3).Runs(1).Text = "daily_usage.csv (~11k rows): Activity logs (Transfers, Connections). Aggregated to per-trial summaries."
$tr3.Paragraphs(6).Runs(1).Text = "Total Samples: 416 complete trials (Filtered from ~500)."
$tr3.Paragraphs(9).Runs(1).Text = "Merged Usage + Subscriptions on 'subscription_id'."
$tr3.Paragraphs(10).Runs(1).Text = "Handled Inactivity (NaN/Zero imputation for missing days)."

# ---------------------------------------------------------------------
# Slide 4: Methodology & Models
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(3).Runs(1).Text = "Tabular Features (157 dims): Sum/Mean/Max/Std of daily activities."
$tr4.Paragraphs(4).Runs(1).Text = "Sequential Data (15 days): Time-series for Deep Learning."
$tr4.Paragraphs(10).Runs(1).Text = "ROC-AUC: Discrimination capability (Primary Metric)."

# ---------------------------------------------------------------------
# Slide 5: Overall Results Comparison
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

# Split the combined ROC-AUC/PR-AUC/Accuracy line into two bullets.
$para5_3 = $tr5.Paragraphs(3)
$para5_3.Runs(1).Text = "ROC-AUC: 0.790 (Best Discrimination)"
$para5_3.Runs(1).InsertAfter($CR + "PR-AUC: 0.835 | Accuracy: 72.3%")

# "Runner Up: GRU" -> "Runner Up: GRU (RNN)" (now paragraph 7 after the split above)
$tr5.Paragraphs(7).Runs(1).Text = "Runner Up: GRU (RNN)"

# Remove the "XGBoost (0.671) - Underperformed LightGBM here." bullet (now paragraph 12)
$tr5.Paragraphs(12).Delete()

# Add a new closing bullet referencing the bar chart.
$tr5.InsertAfter($CR + "   (See Bar Charts ->)")

# ---------------------------------------------------------------------
# Slide 6: Optimization & Training Insights
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Paragraphs(7).Runs(1).Text = "Transformer: Signs of overfitting (Train Loss << Val Loss)."
$tr6.Paragraphs(8).Runs(1).Text = "Takeaway: Deep Learning needs more than 400 samples to outperform Trees."
$tr6.InsertAfter($CR + "   (See Optimization History ->)")

# ---------------------------------------------------------------------
# Slide 7: Feature Importance & Insights
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$tr7.Paragraphs(5).Runs(1).Text = "nb_client_invoices_created_sum: Active usage is the #1 signal."
$tr7.Paragraphs(9).Runs(1).Text = "Targeting: Focus CX on TPEs with low early activity and high churn prob."

# ---------------------------------------------------------------------
# Slide 8: Business Impact & Recommendations
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Paragraphs(5).Runs(1).Text = "Value Calc: 400 * 0.05 * €3k (LTV) = ~€60k/month -> €720k/year."
$tr8.Paragraphs(7).Runs(1).Text = "Day 1-3 (Automated): Nudge users if 'nb_connections' < 2."
$tr8.Paragraphs(8).Runs(1).Text = "Day 7-10 (Human): CX call if Churn Prob > 60%."

# ---------------------------------------------------------------------
# Slide 9: Limitations & Improvements
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Paragraphs(4).Runs(1).Text = "External Factors: No data on economic context or seasonality."
$tr9.Paragraphs(8).Runs(1).Text = "Continuous Training: Retrain monthly to handle data drift."

# ---------------------------------------------------------------------
# Slide 10: Conclusion & Next Steps
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$tr10.InsertAfter($CR + " " + $CR + "Thank You! Questions?")
